$wb = $excel.ActiveWorkbook

# --- Sheet: Bidirectional A ---
$ws = $wb.Worksheets.Item("Bidirectional A")
$ws.Columns.Item(2).ColumnWidth = 22.14
$ws.Range("B2").Value = [double]"0.0006103999985498376"
$ws.Range("B3").Value = [double]"0.0006560000038007274"
$ws.Range("B6").Value = [double]"0.00292205810546875"
$ws.Range("B7").Value = [double]"0.00286712646484375"
$ws.Range("B8").Value = [double]"0.006529808044433594"
$ws.Range("B9").Value = [double]"0.006529808044433594"
$ws.Range("B10").Value = [double]"8.159999561030418e-05"
$ws.Range("B11").Value = [double]"8.369999704882503e-05"
$ws.Range("B12").Value = [double]"7.849099951272365e-05"
$ws.Range("B13").Value = [double]"6.802499992772937e-05"

# --- Sheet: D Lite ---
$ws = $wb.Worksheets.Item("D Lite")
$ws.Range("B2").Value = [double]"0.0003914999979315326"
$ws.Range("B3").Value = [double]"0.002471799998602364"
$ws.Range("B6").Value = [double]"0.0028533935546875"
$ws.Range("B7").Value = [double]"0.0028533935546875"
$ws.Range("B8").Value = [double]"0.07830810546875"
$ws.Range("B9").Value = [double]"0.07811584472656249"
$ws.Range("B10").Value = [double]"0.000622500003373716"
$ws.Range("B11").Value = [double]"0.0008314000006066635"
$ws.Range("B12").Value = [double]"5.738300053053536e-05"
$ws.Range("B13").Value = [double]"2.851799989002757e-05"

# --- Sheet: IDA ---
$ws = $wb.Worksheets.Item("IDA")
$ws.Range("B2").Value = [double]"0.0003122999987681396"
$ws.Range("B3").Value = [double]"0.003447600000072271"
$ws.Range("B6").Value = [double]"0.0028533935546875"
$ws.Range("B7").Value = [double]"0.0028533935546875"
$ws.Range("B8").Value = [double]"0.0027618408203125"
$ws.Range("B9").Value = [double]"0.0026763916015625"
$ws.Range("B10").Value = [double]"0.0001033000007737428"
$ws.Range("B11").Value = [double]"0.003874099995300639"
$ws.Range("B12").Value = [double]"0.001430455000154325"
$ws.Range("B13").Value = [double]"8.190399981685914e-05"

# --- Sheet: SMA ---
$ws = $wb.Worksheets.Item("SMA")
$ws.Range("B2").Value = [double]"0.000745400000596419"
$ws.Range("B3").Value = [double]"0.0002872999975807033"
$ws.Range("B6").Value = [double]"0.002899169921875"
$ws.Range("B7").Value = [double]"0.00289459228515625"
$ws.Range("B8").Value = [double]"0.00261688232421875"
$ws.Range("B9").Value = [double]"0.00261688232421875"
$ws.Range("B10").Value = [double]"8.789999992586672e-05"
$ws.Range("B11").Value = [double]"0.0001325999983237125"
$ws.Range("B12").Value = [double]"0.0001251799998863135"
$ws.Range("B13").Value = [double]"8.085900022706482e-05"

# --- Sheet: RTAA (L=25, M=3) ---
$ws = $wb.Worksheets.Item("RTAA (L=25, M=3)")
$ws.Range("B2").Value = [double]"0.001102899994293693"
$ws.Range("B3").Value = [double]"0.001140100001066457"
$ws.Range("B6").Value = [double]"0.00290679931640625"
$ws.Range("B7").Value = [double]"0.00286407470703125"
$ws.Range("B8").Value = [double]"0.00658416748046875"
$ws.Range("B9").Value = [double]"0.00658416748046875"
$ws.Range("B10").Value = [double]"0.0002367999986745417"
$ws.Range("B11").Value = [double]"0.0005340999996406026"
$ws.Range("B12").Value = [double]"0.0003921499994612532"
$ws.Range("B13").Value = [double]"0.0001180289998592343"
